$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3342080079727115
$ws.Range("C2").Value = 0.8001272719199065
$ws.Range("D2").Value = 1.057948689255649
$ws.Range("E2").Value = 1.02856632710567
$ws.Range("F2").Value = 0.9994115000724915
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = 0.2715251256614695
$ws.Range("C3").Value = 0.5749987921201453
$ws.Range("D3").Value = 0.6574421919879687
$ws.Range("E3").Value = 0.8108280902805284
$ws.Range("F3").Value = 0.7861631834124684
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.1126313932437445
$ws.Range("C4").Value = 0.4526599150000412
$ws.Range("D4").Value = 0.4170638627520434
$ws.Range("E4").Value = 0.645804817845178
$ws.Range("F4").Value = 0.6554781911006021
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.2541178104058748
$ws.Range("C5").Value = 0.4548235677275672
$ws.Range("D5").Value = 0.357578812029527
$ws.Range("E5").Value = 0.5979789394531608
$ws.Range("F5").Value = 0.5590496225693991
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.3625344086844117
$ws.Range("C6").Value = 0.3625344086844117
$ws.Range("D6").Value = 0.1936168257130608
$ws.Range("E6").Value = 0.4400191197130652
$ws.Range("F6").Value = 0.2581229529138574
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.3044989716643398
$ws.Range("C7").Value = 0.3044989716643398
$ws.Range("D7").Value = 0.1348154283086344
$ws.Range("E7").Value = 0.3671722052506622
$ws.Range("F7").Value = 0.212917691553823
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.3319890590706304
$ws.Range("C8").Value = 0.3319890590706304
$ws.Range("D8").Value = 0.1468494055944594
$ws.Range("E8").Value = 0.3832093495655597
$ws.Range("F8").Value = 0.1992119292935162
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.333398904683986
$ws.Range("C9").Value = 0.3399377628279546
$ws.Range("D9").Value = 0.1517941216966595
$ws.Range("E9").Value = 0.3896076509729494
$ws.Range("F9").Value = 0.2105558670468018
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.3614593346982711
$ws.Range("C10").Value = 0.3614593346982711
$ws.Range("D10").Value = 0.1605564165475934
$ws.Range("E10").Value = 0.4006949170473633
$ws.Range("F10").Value = 0.1813668175212442
$ws.Range("G10").Value = 11

$ws.Range("B11").Value = 0.3457552520087644
$ws.Range("C11").Value = 0.3457552520087644
$ws.Range("D11").Value = 0.1504834501066176
$ws.Range("E11").Value = 0.3879219639394211
$ws.Range("F11").Value = 0.1854027322553479
$ws.Range("G11").Value = 10
